$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Made some Changes"
$ws.Range("A3").Value = "Windows 10"
$ws.Range("A5").Value = "SCCM 2012 R2"
$ws.Range("A6").Value = "SCCM CB"
$ws.Range("A7").Value = "SCCM CBB"

$ws.Range("A8").Select()
